$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Bad Drivers" table values that changed for this week's report
$ws.Cells.Item(3, 3).Value2 = 4649

$ws.Cells.Item(5, 2).Value2 = 91
$ws.Cells.Item(5, 3).Value2 = 5137
$ws.Cells.Item(5, 4).Value2 = 97.40000000000001

$ws.Cells.Item(7, 4).Value2 = 98.59999999999999

# The driver "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2" dropped out of the
# Bad Drivers list this week, so remove its row; everything below (including
# the Totals row and the whole Good Drivers table) shifts up by one row.
$ws.Rows.Item(8).Delete()

# Totals row (now row 9 after the deletion) reflects the updated figures
$ws.Cells.Item(9, 2).Value2 = 120
$ws.Cells.Item(9, 3).Value2 = 10652
